# The diff inserts one new data row (row 323) into the "Coliflor" price
# listing and shifts all subsequent rows (old 323-405) down by one, so
# old row 323 becomes new row 324, ..., old row 405 becomes new row 406.
# The sheet's used range grows from A1:R405 to A1:R406.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 323; this shifts rows 323:405 down to 324:406
# and extends the worksheet dimension accordingly.
$ws.Rows.Item(323).Insert()

# Populate the newly inserted row 323 with the new record.
$ws.Cells.Item(323, 1).Value = 10
$ws.Cells.Item(323, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(323, 3).Value = "La Araucanía"
$ws.Cells.Item(323, 4).Value = 44736
$ws.Cells.Item(323, 5).Value = 9
$ws.Cells.Item(323, 6).Value = 100112008
$ws.Cells.Item(323, 7).Value = "Coliflor"
$ws.Cells.Item(323, 8).Value = "Sin especificar"
$ws.Cells.Item(323, 9).Value = "Primera"
$ws.Cells.Item(323, 10).Value = 800
$ws.Cells.Item(323, 11).Value = 1200
$ws.Cells.Item(323, 12).Value = 1200
$ws.Cells.Item(323, 13).Value = 1200
$ws.Cells.Item(323, 14).Value = "`$/unidad"
$ws.Cells.Item(323, 15).Value = "Región Metropolitana"
$ws.Cells.Item(323, 16).Value = 1200
$ws.Cells.Item(323, 17).Value = 1
$ws.Cells.Item(323, 18).Value = "Hortaliza"
